$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Consolidate line items: update SKUs and quantities
$ws.Range("B2").Value = "DS2310BLK-LF"
$ws.Range("B3").Value = "DS2310WMUS-LF"

$ws.Range("C2").Value = 5
$ws.Range("C3").Value = 243
